$d = $word.ActiveDocument

# Simple text replacements where run structure is preserved automatically
# (title heading has no leading empty run; bold/italic runs keep their
# distinct formatting so the adjoining empty run is not merged away)
[void]$d.Content.Find.Execute("Play Hat Trick Hero free - Review of Betsoft Slot Game", $true, $false, $false, $false, $false, $true, 1, $false, "Play Hat Trick Hero Free: Slot Game Review", 2)

[void]$d.Content.Find.Execute("Read our review of Hat Trick Hero, a high volatility slot game by Betsoft with free spin event and potential for big rewards. Play for free now.", $true, $false, $false, $false, $false, $true, 1, $false, "Learn about the features of Hat Trick Hero slot game and play for free.", 2)

# For the "What we like" / "What we don't like" bullet paragraphs, a plain
# Find/Replace on the run's text causes the engine to merge the leading
# empty run (<w:r/>) into the replaced text run when the two runs share the
# same (default) formatting. To keep the original <w:r/><w:r> structure
# intact (matching the source XML exactly), rebuild each of those
# paragraphs via InsertXML instead of a text replace.

function Replace-BulletParagraph($oldText, $newText) {
    $target = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq ($oldText + [char]13)) {
            $target = $p
            break
        }
    }
    if ($target -ne $null) {
        $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>' + $newText + '</w:t></w:r></w:p>'
        [void]$target.Range.InsertXML($xml)
    }
}

Replace-BulletParagraph "Power Shot wild symbol with access to free spins event" "Wild symbol can replace any other symbol"
Replace-BulletParagraph "Collecting Cannon Shot symbols can lead to free spin phase and prize of up to 1000x bet" "Free spins event with expanding wilds"
Replace-BulletParagraph "Free spins event can be extended with additional Cannon Shot symbols" "Chance to win up to 1000 times your bet"
Replace-BulletParagraph "High volatility offers the potential for bigger rewards" "Designed for players who enjoy high volatility"
Replace-BulletParagraph "Payouts may be less frequent due to high volatility" "Payouts may be less frequent"
Replace-BulletParagraph "Goal counter resets at the end of each Hat Trick Hero phase" "Goal counter resets at the end of each phase"
